$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO" - client QUIJIJE MENDOZA GENESIS XIOMARA, row 16
# Column D = "240X80 PORCELANATO": registers a new sale of 413.5
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D16").Value = 413.5
# Row 19 summary counter for column D ("240X80 PORCELANATO") goes from 0 to 1 clients
$wsGrupo.Range("D19").Value = "1 de 17"

# Sheet "VENTA MENSUAL" - same client, row 16, column F = "julio"
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F16").Value = 413.5
# Row 19 totals for "julio" column increases by the same 413.5
$wsMensual.Range("F19").Value = 2057.76

# Sheet "CUMPLIMIENTO MENSUAL" - group totals reflecting the new sale
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Row 2 = "OTROS" group
$wsCumplimiento.Range("D2").Value = 6786.71
$wsCumplimiento.Range("E2").Value = -6786.71
# Row 4 = TOTAL row
$wsCumplimiento.Range("D4").Value = 21006.76
$wsCumplimiento.Range("E4").Value = -7283.419999999999
$wsCumplimiento.Range("F4").Value = 1.530732314436573
